# Integración con Google Drive: se agrega la columna "evidencias" (N) y se
# registran 4 nuevos envíos de formulario (filas 15-18) más un envío en
# proceso (fila 19) que llegó con errores al subir los archivos adjuntos.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Nueva columna "evidencias" (N) con el mismo formato que el resto de
#    encabezados (negrilla, centrado, borde).
# ---------------------------------------------------------------------
$ws.Range("N1").Value = "evidencias"
$ws.Range("N1").Font.Bold = $true
$ws.Range("N1").HorizontalAlignment = -4108
$ws.Range("N1").VerticalAlignment = -4160
$ws.Range("N1").Borders.LineStyle = 1

# Materializar las celdas vacías de la columna N para las filas existentes
# (2-14), igual que el resto de columnas del reporte.
$ws.Range("N2:N14").NumberFormat = "General"

# ---------------------------------------------------------------------
# 2. Fila 14: el "pedido" 23386577 ahora llega como valor numérico (antes
#    se guardaba como texto).
# ---------------------------------------------------------------------
$ws.Range("B14").Value = 23386577

# ---------------------------------------------------------------------
# 3. Nuevos registros del formulario (filas 15 a 19).
# ---------------------------------------------------------------------

# Fila 15
$ws.Range("A15").Value = "2025-10-26 15:34:16"
$ws.Range("B15").Value = 23252885
$ws.Range("C15").Value = "se suben evidencias"
$ws.Range("D15:G15").NumberFormat = "General"
$ws.Range("H15").Value = "CALEB OBED RAMIREZ MUÑOZ"
$ws.Range("I15").Value = "CR 34 E CL 31 -190 (INTERIOR 222 )"
$ws.Range("J15").Value = "Cumplido"
$ws.Range("K15").Value = "VENCIDO"
$ws.Range("L15").Value = 98584619
$ws.Range("M15").Value = "Formulario"
$ws.Range("N15").Value = "Error al subir 23252885_1_20251026_153413.png, Error al subir 23252885_2_20251026_153415.pdf"

# Fila 16
$ws.Range("A16").Value = "2025-10-26 15:50:18"
$ws.Range("B16").Value = 23206334
$ws.Range("C16").Value = "dsasdadsadsa"
$ws.Range("D16:G16").NumberFormat = "General"
$ws.Range("H16").Value = "MARTA ISABEL RODRIGUEZ VELASQUEZ"
$ws.Range("I16").Value = "CL 57 C CR 82 AA -59 (INTERIOR 106 )"
$ws.Range("J16").Value = "Cumplido"
$ws.Range("K16").Value = "VENCIDO"
$ws.Range("L16").Value = 43548242
$ws.Range("M16").Value = "Formulario"
$ws.Range("N16").Value = "Error al subir 23206334_1_20251026_154727.png, Error al subir 23206334_2_20251026_155016.pdf"

# Fila 17
$ws.Range("A17").Value = "2025-10-26 15:51:45"
$ws.Range("B17").Value = 23252901
$ws.Range("C17").Value = "ddasddasdadas"
$ws.Range("D17:G17").NumberFormat = "General"
$ws.Range("H17").Value = "CALEB OBED RAMIREZ MUÑOZ"
$ws.Range("I17").Value = "CR 34 E CL 31 -190 (INTERIOR 127 )"
$ws.Range("J17").Value = "Cumplido"
$ws.Range("K17").Value = "VENCIDO"
$ws.Range("L17").Value = 98584619
$ws.Range("M17").Value = "Formulario"
$ws.Range("N17").Value = "Error al subir 23252901_1_20251026_155141.png, Error al subir 23252901_2_20251026_155143.pdf"

# Fila 18
$ws.Range("A18").Value = "2025-10-26 15:56:43"
$ws.Range("B18").Value = 23263902
$ws.Range("C18").Value = "suba suba por fa"
$ws.Range("D18:G18").NumberFormat = "General"
$ws.Range("H18").Value = "GLORIA PATRICIA ZULUAGA GOMEZ"
$ws.Range("I18").Value = "CR 144 CL 68 -143"
$ws.Range("J18").Value = "Cumplido"
$ws.Range("K18").Value = "VENCIDO"
$ws.Range("L18").Value = 43588429
$ws.Range("M18").Value = "Formulario"
$ws.Range("N18").Value = "Error al subir 23263902_1_20251026_155639.pdf, Error al subir 23263902_2_20251026_155641.png"

# Fila 19 (pedido llega como texto, igual que en el archivo original)
$ws.Range("A19").Value = "2025-10-26 15:57:31"
$ws.Range("B19").Value = "'23305567"
$ws.Range("C19").Value = "sdsaaddas"
$ws.Range("D19:G19").NumberFormat = "General"
$ws.Range("H19").Value = "LUIS MARIANO UPEGUI FERNANDEZ"
$ws.Range("I19").Value = "CR 29 CL 6 -24 (INTERIOR 5020 )"
$ws.Range("J19").Value = "Cumplido"
$ws.Range("K19").Value = "VENCIDO"
$ws.Range("L19").Value = 71639305
$ws.Range("M19").Value = "Formulario"
$ws.Range("N19").Value = "Error al subir 23305567_1_20251026_155727.png, Error al subir 23305567_2_20251026_155729.pdf"
